$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.042152351287598
$ws.Range("C2").Value = 0.2477577134807518
$ws.Range("D2").Value = 0.07932116302885106
$ws.Range("E2").Value = 0.4213043912050125
$ws.Range("G2").Value = 0.00239256667934717
$ws.Range("N2").Value = 0.8858854290019202
$ws.Range("O2").Value = 2.36159559119335
$ws.Range("B3").Value = 0.9216377364107871
$ws.Range("C3").Value = 0.21634972604312
$ws.Range("D3").Value = 0.07182853882660822
$ws.Range("E3").Value = 0.3673556333177004
$ws.Range("G3").Value = 0.002396467511749116
$ws.Range("N3").Value = 0.8986418128053337
$ws.Range("O3").Value = 2.311976073344738
$ws.Range("B4").Value = 0.8477359112253566
$ws.Range("C4").Value = 0.1970005991078381
$ws.Range("D4").Value = 0.06726620138621797
$ws.Range("E4").Value = 0.3343431044413165
$ws.Range("G4").Value = 0.002398987038450047
$ws.Range("N4").Value = 0.9069359562920027
$ws.Range("O4").Value = 2.283618760817149
$ws.Range("B5").Value = 0.8176441028412
$ws.Range("C5").Value = 0.1890994255425085
$ws.Range("D5").Value = 0.06541651123586689
$ws.Range("E5").Value = 0.3209161034681074
$ws.Range("G5").Value = 0.002400045155102981
$ws.Range("N5").Value = 0.9104318575920836
$ws.Range("O5").Value = 2.272589650484605
$ws.Range("B6").Value = 0.8126488178158979
$ws.Range("C6").Value = 0.18778645969374
$ws.Range("D6").Value = 0.06510994252035118
$ws.Range("E6").Value = 0.3186880546108881
$ws.Range("G6").Value = 0.00240022275360274
$ws.Range("N6").Value = 0.9110193502662511
$ws.Range("O6").Value = 2.270789981166274
$ws.Range("B7").Value = 0.8473299862673684
$ws.Range("C7").Value = 0.1968941069328025
$ws.Range("D7").Value = 0.06724121748173673
$ws.Range("E7").Value = 0.3341619215214706
$ws.Range("G7").Value = 0.002399001181207333
$ws.Range("N7").Value = 0.9069826338439597
$ws.Range("O7").Value = 2.28346789055999
$ws.Range("B8").Value = 1.000579220376949
$ws.Range("C8").Value = 0.2369415970946136
$ws.Range("D8").Value = 0.07672972408518319
$ws.Range("E8").Value = 0.4026781333234339
$ws.Range("G8").Value = 0.002393885932271027
$ws.Range("N8").Value = 0.8901879366963463
$ws.Range("O8").Value = 2.344046765480954
$ws.Range("B9").Value = 1.301868333419861
$ws.Range("C9").Value = 0.3149682249249963
$ws.Range("D9").Value = 0.09564487307289937
$ws.Range("E9").Value = 0.5380411872147874
$ws.Range("G9").Value = 0.002384837066387848
$ws.Range("N9").Value = 0.8609220568304394
$ws.Range("O9").Value = 2.479753539211174
$ws.Range("B10").Value = 1.523740703338603
$ws.Range("C10").Value = 0.3719989580362721
$ws.Range("D10").Value = 0.1097386695375207
$ws.Range("E10").Value = 0.6382692049963339
$ws.Range("G10").Value = 0.002378780664825875
$ws.Range("N10").Value = 0.8416639627964813
$ws.Range("O10").Value = 2.590025063619919
$ws.Range("B11").Value = 1.624799985612867
$ws.Range("C11").Value = 0.3978828608658205
$ws.Range("D11").Value = 0.1161950717264801
$ws.Range("E11").Value = 0.6840708183224251
$ws.Range("G11").Value = 0.002376152474220963
$ws.Range("N11").Value = 0.8333917987594148
$ws.Range("O11").Value = 2.642542055918113
$ws.Range("B12").Value = 1.663087465220087
$ws.Range("C12").Value = 0.407675991914914
$ws.Range("D12").Value = 0.1186465576222275
$ws.Range("E12").Value = 0.7014473482350496
$ws.Range("G12").Value = 0.002375175381049807
$ws.Range("N12").Value = 0.8303297557048595
$ws.Range("O12").Value = 2.662771758954023
$ws.Range("B13").Value = 1.654840746502089
$ws.Range("C13").Value = 0.4055672451134456
$ws.Range("D13").Value = 0.118118291625791
$ws.Range("E13").Value = 0.6977035159659266
$ws.Range("G13").Value = 0.002375385010077226
$ws.Range("N13").Value = 0.8309860849446906
$ws.Range("O13").Value = 2.658399623855587
$ws.Range("B14").Value = 1.627949547980677
$ws.Range("C14").Value = 0.3986887185857881
$ws.Range("D14").Value = 0.1163966243283312
$ws.Range("E14").Value = 0.6854997302826717
$ws.Range("G14").Value = 0.00237607172495862
$ws.Range("N14").Value = 0.8331384695677286
$ws.Range("O14").Value = 2.644199474378013
$ws.Range("B15").Value = 1.611480326657215
$ws.Range("C15").Value = 0.394474311411102
$ws.Range("D15").Value = 0.1153429143367646
$ws.Range("E15").Value = 0.6780288642481906
$ws.Range("G15").Value = 0.00237649471777737
$ws.Range("N15").Value = 0.8344660477778802
$ws.Range("O15").Value = 2.63554622132267
$ws.Range("B16").Value = 1.517138830021679
$ws.Range("C16").Value = 0.3703061866565918
$ws.Range("D16").Value = 0.1093176443307726
$ws.Range("E16").Value = 0.6352803590481244
$ws.Range("G16").Value = 0.002378954968118493
$ws.Range("N16").Value = 0.8422144144945136
$ws.Range("O16").Value = 2.586640676848447
$ws.Range("B17").Value = 1.459296302616622
$ws.Range("C17").Value = 0.3554645688490723
$ws.Range("D17").Value = 0.1056329534862641
$ws.Range("E17").Value = 0.6091105125323395
$ws.Range("G17").Value = 0.002380496680864501
$ws.Range("N17").Value = 0.8470930336244749
$ws.Range("O17").Value = 2.557244557156537
$ws.Range("B18").Value = 1.426038894238843
$ws.Range("C18").Value = 0.3469224442941083
$ws.Range("D18").Value = 0.1035178519105955
$ws.Range("E18").Value = 0.5940776694852588
$ws.Range("G18").Value = 0.002381395383206441
$ws.Range("N18").Value = 0.8499450475608228
$ws.Range("O18").Value = 2.540557997430085
$ws.Range("B19").Value = 1.414780581282855
$ws.Range("C19").Value = 0.3440292616461988
$ws.Range("D19").Value = 0.102802439028963
$ws.Range("E19").Value = 0.5889910527613438
$ws.Range("G19").Value = 0.002381701724035623
$ws.Range("N19").Value = 0.8509185781239523
$ws.Range("O19").Value = 2.534946099050842
$ws.Range("B20").Value = 1.465452487475204
$ws.Range("C20").Value = 0.357045065630416
$ws.Range("D20").Value = 0.1060247560255902
$ws.Range("E20").Value = 0.6118943146406224
$ws.Range("G20").Value = 0.002380331326716912
$ws.Range("N20").Value = 0.8465689380096393
$ws.Range("O20").Value = 2.560350889465951
$ws.Range("B21").Value = 1.635847641582302
$ws.Range("C21").Value = 0.4007093400332451
$ws.Range("D21").Value = 0.1169021401369719
$ws.Range("E21").Value = 0.6890833760429302
$ws.Range("G21").Value = 0.002375869528450223
$ws.Range("N21").Value = 0.8325043484405654
$ws.Range("O21").Value = 2.648361069449265
$ws.Range("B22").Value = 1.747318732050132
$ws.Range("C22").Value = 0.4291967054154497
$ws.Range("D22").Value = 0.1240495810455542
$ws.Range("E22").Value = 0.739721298735617
$ws.Range("G22").Value = 0.002373059206490188
$ws.Range("N22").Value = 0.8237230417123627
$ws.Range("O22").Value = 2.707880044962565
$ws.Range("B23").Value = 1.68781460614224
$ws.Range("C23").Value = 0.4139970012097365
$ws.Range("D23").Value = 0.120231304384987
$ws.Range("E23").Value = 0.7126765847421268
$ws.Range("G23").Value = 0.002374549488305271
$ws.Range("N23").Value = 0.8283721387044309
$ws.Range("O23").Value = 2.675929329946939
$ws.Range("B24").Value = 1.462669285555251
$ws.Range("C24").Value = 0.3563305525744909
$ws.Range("D24").Value = 0.1058476119019502
$ws.Range("E24").Value = 0.6106357188408538
$ws.Range("G24").Value = 0.002380406044883362
$ws.Range("N24").Value = 0.8468057345283668
$ws.Range("O24").Value = 2.558945851403905
$ws.Range("B25").Value = 1.22027465209726
$ws.Range("C25").Value = 0.2939128168450793
$ws.Range("D25").Value = 0.09049390800328183
$ws.Range("E25").Value = 0.5012975497999435
$ws.Range("G25").Value = 0.002387180599359666
$ws.Range("N25").Value = 0.8684458098405834
$ws.Range("O25").Value = 2.441204220137791
